$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12: 2015-05-28 -> "Learn Requirejs, Demo project: load more comment, delete comment"
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats - reuse existing date style
$ws.Range("A12").Value2 = 42152
$ws.Range("B12").Value = "Learn Requirejs, Demo project: load more comment, delete comment"

# New row 13: 2015-05-29 -> "Spring Security (the demo project has CSRF error in delete comment feature)"
$ws.Range("A11").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats - reuse existing date style
$ws.Range("A13").Value2 = 42153
$ws.Range("B13").Value = "Spring Security (the demo project has CSRF error in delete comment feature)"

$excel.CutCopyMode = $false

# Update selection to match the new active cell after edits
$ws.Range("B13").Select()
